$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price ticker refresh: update Price (D) and Volume(1h) (E) columns
# for rows 2-51 to match the latest scrape. D-column values that parse as plain
# numbers (no thousands separators) are forced to remain text, matching the
# original inline-string cell contents, then formatting is reset to the default
# (unstyled) cell style so no stray number-format style sticks around.

$ws.Range("D2").Value = '27.165.60'
$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").Value = '1.677.90'
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  -0.04%  '

$d = $ws.Range("D5")
$d.NumberFormat = '@'
$d.Value = '214.38'
$d.ClearFormats()
$ws.Range("E5").Value = '  -0.50%  '

$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("E7").Value = '  -0.08%  '

$d = $ws.Range("D8")
$d.NumberFormat = '@'
$d.Value = '22.97'
$d.ClearFormats()
$ws.Range("E8").Value = '  +7.87%  '

$ws.Range("E9").Value = '  +3.35%  '

$ws.Range("E10").Value = '  +0.01%  '

$d = $ws.Range("D11")
$d.NumberFormat = '@'
$d.Value = '0.0890'
$d.ClearFormats()
$ws.Range("E11").Value = '  +0.21%  '

$ws.Range("D12").Value = '1.914.81'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = '1.684.44'
$ws.Range("E13").Value = '  +0.04%  '

$ws.Range("E14").Value = '  +2.17%  '

$d = $ws.Range("D15")
$d.NumberFormat = '@'
$d.Value = '0.560'
$d.ClearFormats()
$ws.Range("E15").Value = '  +4.87%  '

$d = $ws.Range("D16")
$d.NumberFormat = '@'
$d.Value = '66.63'
$d.ClearFormats()
$ws.Range("E16").Value = '  +0.50%  '

$ws.Range("D17").Value = '27.147.72'
$ws.Range("E17").Value = '  +0.44%  '

$d = $ws.Range("D18")
$d.NumberFormat = '@'
$d.Value = '235.55'
$d.ClearFormats()
$ws.Range("E18").Value = '  +0.37%  '

$ws.Range("D19").Value = '0.0₃0742'
$ws.Range("E19").Value = '  +1.14%  '

$d = $ws.Range("D20")
$d.NumberFormat = '@'
$d.Value = '7.84'
$d.ClearFormats()
$ws.Range("E20").Value = '  -3.61%  '

$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("E22").Value = '  +2.31%  '

$ws.Range("E23").Value = '  +3.44%  '

$d = $ws.Range("D24")
$d.NumberFormat = '@'
$d.Value = '2.10'
$d.ClearFormats()
$ws.Range("E24").Value = '  -0.84%  '

$d = $ws.Range("D25")
$d.NumberFormat = '@'
$d.Value = '147.64'
$d.ClearFormats()
$ws.Range("E25").Value = '  +0.28%  '

$ws.Range("E26").Value = '  +3.08%  '

$d = $ws.Range("D27")
$d.NumberFormat = '@'
$d.Value = '16.43'
$d.ClearFormats()
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("E28").Value = '  +0.46%  '

$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("E30").Value = '  +0.60%  '

$ws.Range("E31").Value = '  -0.01%  '

$ws.Range("E32").Value = '  +0.18%  '

$ws.Range("D33").Value = '1.537.62'
$ws.Range("E33").Value = '  -0.21%  '

$ws.Range("E34").Value = '  +2.65%  '

$ws.Range("E35").Value = '  -2.98%  '

$d = $ws.Range("D36")
$d.NumberFormat = '@'
$d.Value = '0.609'
$d.ClearFormats()
$ws.Range("E36").Value = '  +4.50%  '

$ws.Range("E37").Value = '  +4.34%  '

$ws.Range("E38").Value = '  -0.31%  '

$ws.Range("E39").Value = '  -0.65%  '

$ws.Range("E40").Value = '  +2.14%  '

$d = $ws.Range("D41")
$d.NumberFormat = '@'
$d.Value = '69.86'
$d.ClearFormats()
$ws.Range("E41").Value = '  +3.07%  '

$d = $ws.Range("D42")
$d.NumberFormat = '@'
$d.Value = '5.79'
$d.ClearFormats()
$ws.Range("E42").Value = '  +4.68%  '

$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("E44").Value = '  +0.25%  '

$ws.Range("D45").Value = '1.823.93'
$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("E46").Value = '  +0.35%  '

$ws.Range("E47").Value = '  +7.61%  '

$d = $ws.Range("D48")
$d.NumberFormat = '@'
$d.Value = '90.13'
$d.ClearFormats()
$ws.Range("E48").Value = '  -0.19%  '

$ws.Range("E49").Value = '  +2.78%  '

$d = $ws.Range("D50")
$d.NumberFormat = '@'
$d.Value = '8.22'
$d.ClearFormats()
$ws.Range("E50").Value = '  +3.17%  '

$ws.Range("E51").Value = '  +1.29%  '
